# Update "Horario Final" column (D2:D10) on the "impressoras" sheet to
# 23:59:59 for every printer row, then leave the selection on E6 as the
# last thing the author touched while finishing up / reviewing logs.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("impressoras")

$novoHorario = 0.99998842592592596 # 23:59:59 as an Excel serial time fraction

for ($row = 2; $row -le 10; $row++) {
    $ws.Cells.Item($row, 4).Value = $novoHorario
}

$ws.Activate()
$ws.Range("E6").Select()
